$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 233, shifting existing rows
# (old 233..256) down to (234..257).
$ws.Rows.Item(233).Insert()

# Populate the newly inserted row 233 with the new weekly price record.
$ws.Range("A233").Value = 8
$ws.Range("B233").Value = "Terminal La Palmera de La Serena"
$ws.Range("C233").Value = "Coquimbo"
$ws.Range("D233").Value = [DateTime]"2023-06-29"
$ws.Range("E233").Value = 4
$ws.Range("F233").Value = 100112001
$ws.Range("G233").Value = "Berenjena"
$ws.Range("H233").Value = "Sin especificar"
$ws.Range("I233").Value = "Primera"
$ws.Range("J233").Value = 240
$ws.Range("K233").Value = 8000
$ws.Range("L233").Value = 9000
$ws.Range("M233").Value = 8500
$ws.Range("N233").Value = "$/caja 50 unidades"
$ws.Range("O233").Value = "Región de Arica y Parinacota"
$ws.Range("P233").Value = 170
$ws.Range("Q233").Value = 50
$ws.Range("R233").Value = "Hortaliza"
